$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename category "Other writing" -> "Other work" wherever it occurs
# (the category column, A, for the "Other writing" entries)
$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Text -eq "Other writing") {
        $cell.Value = "Other work"
    }
}

# Update the active selection to A8
$ws.Range("A8").Select()
